$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the frozen-pane view (no more split/pane + custom selection) ---
$excel.ActiveWindow.FreezePanes = $false

# --- Row 1 gets a shorter, explicit custom height instead of the old autosized one ---
$ws.Rows("1").RowHeight = 55.5

# --- Bring in a new "2023" column (H) that mirrors column G's look ---
# Copy the formatting of the existing 2022 column (G4:G34) onto H4:H34 so the
# new column inherits the same number formats / fonts / borders per row.
$ws.Range("G4:G34").Copy()
$ws.Range("H4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header year ---
$ws.Range("H4").Value = 2023

# --- Kyrgyz Republic block (rows 6-19) ---
$ws.Range("H6").Value = 4.684701603429354
$ws.Range("H8").Value = 9.1619177226137172
$ws.Range("H9").Value = 2.1026629736131146
$ws.Range("H11").Value = 1.6883952033845095
$ws.Range("H12").Value = 5.4628226647491029
$ws.Range("H13").Value = 3.5472154122985047
$ws.Range("H14").Value = 5.0893891878983322
$ws.Range("H15").Value = 0.47898402261907741
$ws.Range("H16").Value = 6.7601470318978496
$ws.Range("H17").Value = 3.575054136227692
$ws.Range("H18").Value = 5.6091986107002105
$ws.Range("H19").Value = 23.492624647686782

# --- Second block / territory rows (21-34) ---
$ws.Range("H21").Value = 1.6110270246454947
$ws.Range("H23").Value = 3.5911978326282963
$ws.Range("H24").Value = 0.46905009439682216
$ws.Range("H26").Value = 0.17806447072843995
$ws.Range("H27").Value = 2.6322222213217694
$ws.Range("H28").Value = 1.8451824926932785
$ws.Range("H29").Value = 0.59506627968271797
$ws.Range("H30").Value = 0.028565820593621703
$ws.Range("H31").Value = 0.17527311250618646
$ws.Range("H32").Value = 1.6783444058510675
$ws.Range("H33").Value = 0.45739674737486102
$ws.Range("H34").Value = 11.499084267608914
